$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in day numbers 1..31 into A2:A32
for ($i = 1; $i -le 31; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $i
}

# Update the view: select A2:A32 (with A2 as the active/top-left cell of the
# selection, which also becomes the active cell), scroll so A2 is the
# top-left visible cell, and zoom to 115%.
$excel.Goto($ws.Range("A2:A32"))
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$win.Zoom = 115
